$wb = $excel.ActiveWorkbook

# --- 1. Clean up "ODI Batting Extra": rows whose batting stats are all blank
#        should have no B/C/D/E cells at all (only A + F remain).
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$blankBattingRows = @(4, 6, 14, 16, 17)
foreach ($r in $blankBattingRows) {
    $battingExtra.Range("B$r`:E$r").ClearContents()
}
# Row 10 only has its PERCENT_RUNS_OF_TOTAL (E) blank; B/C/D keep real values.
$battingExtra.Range("E10").ClearContents()

# --- 2. Add the new "ODI Bowling Extra" sheet, positioned right after
#        "ODI Batting Extra" (matches the workbook.xml <sheet> ordering in the diff).
$bowlingExtra = $wb.Worksheets.Add()
$bowlingExtra.Name = "ODI Bowling Extra"
# NOTE: sheet object handles captured before a structural change (Add/Move/Delete)
# can resolve to the wrong sheet afterwards, since they track a position rather than
# a stable identity in this host. Always re-fetch by name right before using a handle
# that survived a structural mutation.
$battingExtraFresh = $wb.Worksheets.Item("ODI Batting Extra")
$bowlingExtra.Move($null, $battingExtraFresh)
$bowlingExtra = $wb.Worksheets.Item("ODI Bowling Extra")

# --- 3. Populate header row
$bowlingExtra.Cells.Item(1, 1).Value = "MATCH_CODE"
$bowlingExtra.Cells.Item(1, 2).Value = "MAIDEN_OVERS"
$bowlingExtra.Cells.Item(1, 3).Value = "PERCENT_WICKETS_OF_ALL"

# --- 4. Populate data rows (MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL)
#        Missing MAIDEN_OVERS/PERCENT_WICKETS_OF_ALL values are intentionally left blank.
$bowlingExtra.Cells.Item(2, 1).Value = "3215"
$bowlingExtra.Cells.Item(2, 2).Value = "0"
$bowlingExtra.Cells.Item(2, 3).Value = "10.00%"
$bowlingExtra.Cells.Item(3, 1).Value = "3217"
$bowlingExtra.Cells.Item(4, 1).Value = "3219"
$bowlingExtra.Cells.Item(4, 2).Value = "0"
$bowlingExtra.Cells.Item(5, 1).Value = "3222"
$bowlingExtra.Cells.Item(5, 2).Value = "0"
$bowlingExtra.Cells.Item(6, 1).Value = "3353"
$bowlingExtra.Cells.Item(6, 2).Value = "0"
$bowlingExtra.Cells.Item(7, 1).Value = "3359"
$bowlingExtra.Cells.Item(7, 2).Value = "0"
$bowlingExtra.Cells.Item(8, 1).Value = "3367"
$bowlingExtra.Cells.Item(8, 2).Value = "0"
$bowlingExtra.Cells.Item(8, 3).Value = "10.00%"
$bowlingExtra.Cells.Item(9, 1).Value = "3373"
$bowlingExtra.Cells.Item(10, 1).Value = "3375"
$bowlingExtra.Cells.Item(10, 2).Value = "0"
$bowlingExtra.Cells.Item(11, 1).Value = "3380"
$bowlingExtra.Cells.Item(11, 2).Value = "0"
$bowlingExtra.Cells.Item(11, 3).Value = "10.00%"
$bowlingExtra.Cells.Item(12, 1).Value = "3395"
$bowlingExtra.Cells.Item(12, 2).Value = "0"
$bowlingExtra.Cells.Item(13, 1).Value = "3397"
$bowlingExtra.Cells.Item(13, 2).Value = "0"
$bowlingExtra.Cells.Item(14, 1).Value = "3399"
$bowlingExtra.Cells.Item(15, 1).Value = "3433"
$bowlingExtra.Cells.Item(15, 2).Value = "0"
$bowlingExtra.Cells.Item(16, 1).Value = "3434"
$bowlingExtra.Cells.Item(17, 1).Value = "3435"
$bowlingExtra.Cells.Item(17, 2).Value = "0"
$bowlingExtra.Cells.Item(18, 1).Value = "3580"
$bowlingExtra.Cells.Item(19, 1).Value = "3621"
$bowlingExtra.Cells.Item(20, 1).Value = "3781"
$bowlingExtra.Cells.Item(20, 2).Value = "0"
$bowlingExtra.Cells.Item(21, 1).Value = "3874"
$bowlingExtra.Cells.Item(21, 2).Value = "0"

